$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 16. This splits the old row 15 (which held
#    both "Nepal contributes..." and "Hence, I strongly call upon..." in a
#    single cell) into two rows, and pushes the former rows 16-21 down to
#    17-22.
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Insert()

# ---------------------------------------------------------------------------
# 2. Row 15 now only keeps the first sentence, and is no longer a "yes"
#    (annotated) row - its analysis columns move to the new row 16.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Nepal contributes less than 0.1 percent of global green house gas emission, yet we suffer the most. "
$ws.Range("B15").Value = "no"
$ws.Range("C15:H15").ClearContents()

# ---------------------------------------------------------------------------
# 3. New row 16 carries the second sentence plus its own annotation.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Hence, I strongly call upon the countries with higher emissions to take leadership in reducing their emissions to secure a safer planet."
$ws.Range("B16").Value = "yes"
$ws.Range("C16").Value = "moral responsibility"
$ws.Range("D16").Value = "emissions"
$ws.Range("E16").Value = "other(countries with higer emissions)"
$ws.Range("F16").Value = "n.a. "
$ws.Range("G16").Value = "egalitarian"
$ws.Range("H16").Value = "Indication of foundational belief that others that emit more have the responsibility to take on action. "

# ---------------------------------------------------------------------------
# 4. Row 10 - principle/topic relabelled.
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = "policy target, temperature"
$ws.Range("G10").Value = "egalitarian"

# ---------------------------------------------------------------------------
# 5. Row 11 - explanation reworded.
# ---------------------------------------------------------------------------
$ws.Range("H11").Value = "Urge to increase the support of the most vulnerable. "

# ---------------------------------------------------------------------------
# 6. Row 21 (previously row 20, shifted down by the insert) - principle and
#    explanation reworded.
# ---------------------------------------------------------------------------
$ws.Range("G21").Value = "egalitarian"
$ws.Range("H21").Value = "Call for a fair agreement, indicating egalitarian"

# ---------------------------------------------------------------------------
# 7. Row heights - match the re-flowed layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 100.8
$ws.Rows.Item(3).RowHeight = 83.4
$ws.Rows.Item(4).RowHeight = 100.8
$ws.Rows.Item(5).RowHeight = 72
$ws.Rows.Item(6).RowHeight = 57.6
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 57.6
$ws.Rows.Item(11).RowHeight = 92.4
$ws.Rows.Item(12).RowHeight = 57.6
$ws.Rows.Item(13).RowHeight = 115.2
$ws.Rows.Item(14).RowHeight = 57.6
$ws.Rows.Item(15).RowHeight = 43.2
$ws.Rows.Item(16).RowHeight = 69.6
$ws.Rows.Item(17).RowHeight = 57.6
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 109.2
$ws.Rows.Item(20).RowHeight = 72
$ws.Rows.Item(21).RowHeight = 123.6

# ---------------------------------------------------------------------------
# 8. Column A width tweaks very slightly to match new font metrics.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.83

# ---------------------------------------------------------------------------
# 9. Selection / scroll position.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C10").Select()

# ---------------------------------------------------------------------------
# 10. Touch the bold font so its family attribute is (re)written.
# ---------------------------------------------------------------------------
$ws.Range("A1").Font.Name = "Calibri"
